# Apply cryptos price/ranking update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "66.222.97"
$ws.Cells.Item(2, 5).Value = "  +1.89%  "
$ws.Cells.Item(3, 4).Value = "3.417.43"
$ws.Cells.Item(3, 5).Value = "  +0.80%  "
$ws.Cells.Item(4, 5).Value = "  -0.01%  "
$ws.Cells.Item(5, 4).Value = "'566.69"
$ws.Cells.Item(5, 5).Value = "  +1.37%  "
$ws.Cells.Item(6, 4).Value = "'178.26"
$ws.Cells.Item(6, 5).Value = "  +2.73%  "
$ws.Cells.Item(7, 5).Value = "  +1.38%  "
$ws.Cells.Item(8, 4).Value = "3.411.64"
$ws.Cells.Item(8, 5).Value = "  +0.87%  "
$ws.Cells.Item(9, 5).Value = "  +0.02%  "
$ws.Cells.Item(10, 4).Value = "'0.176"
$ws.Cells.Item(10, 5).Value = "  +4.32%  "
$ws.Cells.Item(11, 4).Value = "'0.639"
$ws.Cells.Item(11, 5).Value = "  +1.32%  "
$ws.Cells.Item(12, 4).Value = "'54.49"
$ws.Cells.Item(12, 5).Value = "  +0.06%  "
$ws.Cells.Item(13, 4).Value = "'0.0000280"
$ws.Cells.Item(13, 5).Value = "  +0.46%  "
$ws.Cells.Item(14, 5).Value = "  +2.25%  "
$ws.Cells.Item(15, 4).Value = "3.951.83"
$ws.Cells.Item(15, 5).Value = "  +0.59%  "
$ws.Cells.Item(16, 2).Value = "WrappedEther"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(16, 4).Value = "3.418.18"
$ws.Cells.Item(16, 5).Value = "  +0.98%  "
$ws.Cells.Item(17, 2).Value = "Chainlink"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(17, 4).Value = "'18.33"
$ws.Cells.Item(17, 5).Value = "  +0.19%  "
$ws.Cells.Item(18, 2).Value = "TRON"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Cells.Item(18, 4).Value = "'0.120"
$ws.Cells.Item(18, 5).Value = "  +0.91%  "
$ws.Cells.Item(19, 4).Value = "66.111.88"
$ws.Cells.Item(19, 5).Value = "  +1.64%  "
$ws.Cells.Item(20, 4).Value = "'11.97"
$ws.Cells.Item(20, 5).Value = "  +1.34%  "
$ws.Cells.Item(21, 4).Value = "'1.00"
$ws.Cells.Item(21, 5).Value = "  +1.19%  "
$ws.Cells.Item(22, 4).Value = "'466.84"
$ws.Cells.Item(22, 5).Value = "  -0.97%  "
$ws.Cells.Item(23, 4).Value = "'4.95"
$ws.Cells.Item(23, 5).Value = "  -0.30%  "
$ws.Cells.Item(24, 4).Value = "'14.74"
$ws.Cells.Item(24, 5).Value = "  +8.61%  "
$ws.Cells.Item(25, 5).Value = "  +0.31%  "
$ws.Cells.Item(26, 4).Value = "'90.00"
$ws.Cells.Item(26, 5).Value = "  +3.48%  "
$ws.Cells.Item(27, 5).Value = "  +1.31%  "
$ws.Cells.Item(28, 4).Value = "'10.77"
$ws.Cells.Item(28, 5).Value = "  -0.08%  "
$ws.Cells.Item(29, 5).Value = "  +0.79%  "
$ws.Cells.Item(30, 5).Value = "  +1.57%  "
$ws.Cells.Item(31, 4).Value = "'6.76"
$ws.Cells.Item(31, 5).Value = "  +0.94%  "
$ws.Cells.Item(32, 4).Value = "'11.56"
$ws.Cells.Item(32, 5).Value = "  +0.52%  "
$ws.Cells.Item(33, 4).Value = "'583.33"
$ws.Cells.Item(33, 5).Value = "  +1.87%  "
$ws.Cells.Item(34, 5).Value = "  +1.78%  "
$ws.Cells.Item(35, 5).Value = "  +0.76%  "
$ws.Cells.Item(36, 5).Value = "  -0.02%  "
$ws.Cells.Item(37, 4).Value = "'0.145"
$ws.Cells.Item(37, 5).Value = "  +3.68%  "
$ws.Cells.Item(38, 4).Value = "'3.57"
$ws.Cells.Item(38, 5).Value = "  +0.24%  "
$ws.Cells.Item(39, 4).Value = "'36.44"
$ws.Cells.Item(39, 5).Value = "  +2.02%  "
$ws.Cells.Item(40, 2).Value = "PEPE"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Cells.Item(40, 4).Value = "0.0₃0769"
$ws.Cells.Item(40, 5).Value = "  +2.56%  "
$ws.Cells.Item(41, 2).Value = "TheGraph"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Cells.Item(41, 4).Value = "'0.382"
$ws.Cells.Item(41, 5).Value = "  +3.33%  "
$ws.Cells.Item(42, 4).Value = "3.128.85"
$ws.Cells.Item(42, 5).Value = "  +1.07%  "
$ws.Cells.Item(43, 4).Value = "'2.94"
$ws.Cells.Item(43, 5).Value = "  +3.31%  "
$ws.Cells.Item(45, 5).Value = "  +1.64%  "
$ws.Cells.Item(47, 5).Value = "  -0.04%  "
$ws.Cells.Item(48, 4).Value = "'2.63"
$ws.Cells.Item(48, 5).Value = "  +13.65%  "
$ws.Cells.Item(49, 5).Value = "  +0.21%  "
$ws.Cells.Item(50, 4).Value = "'142.09"
$ws.Cells.Item(50, 5).Value = "  +1.93%  "
$ws.Cells.Item(51, 4).Value = "'8.56"
$ws.Cells.Item(51, 5).Value = "  +2.80%  "
